$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 53: copy style (incl. date number format) from A52 so the new date cell formats correctly
$ws.Range("A52").Copy() | Out-Null
$ws.Range("A53").PasteSpecial(-4122) | Out-Null

# row 2
$ws.Range("A2").Value2 = 39400
$ws.Range("B2").Value2 = 2007
$ws.Range("C2").Value2 = 2.070003986395053
$ws.Range("D2").Value2 = 2008
$ws.Range("E2").Value2 = -0.3549868696899106

# row 3
$ws.Range("A3").Value2 = 39583
$ws.Range("B3").Value2 = 2008
$ws.Range("C3").Value2 = 0.5453776865001148
$ws.Range("D3").Value2 = 2009
$ws.Range("E3").Value2 = -0.7976031984000098

# row 4
$ws.Range("A4").Value2 = 39765
$ws.Range("B4").Value2 = 2008
$ws.Range("C4").Value2 = 0.517569958955022
$ws.Range("D4").Value2 = 2009
$ws.Range("E4").Value2 = -5.168396053267498

# row 5
$ws.Range("A5").Value2 = 39948
$ws.Range("B5").Value2 = 2009
$ws.Range("C5").Value2 = -6.170514117037273
$ws.Range("D5").Value2 = 2010
$ws.Range("E5").Value2 = -8.396348489509153

# row 6
$ws.Range("A6").Value2 = 40130
$ws.Range("B6").Value2 = 2009
$ws.Range("C6").Value2 = -3.956152295564896
$ws.Range("D6").Value2 = 2010
$ws.Range("E6").Value2 = -1.314964327391877

# row 7
$ws.Range("A7").Value2 = 40310
$ws.Range("B7").Value2 = 2010
$ws.Range("C7").Value2 = -0.2290082001396909
$ws.Range("D7").Value2 = 2011
$ws.Range("E7").Value2 = -4.327930935900004

# row 8
$ws.Range("A8").Value2 = 40494
$ws.Range("B8").Value2 = 2010
$ws.Range("C8").Value2 = 1.234995474941392
$ws.Range("D8").Value2 = 2011
$ws.Range("E8").Value2 = 1.001424185348321

# row 9
$ws.Range("A9").Value2 = 40676
$ws.Range("B9").Value2 = 2011
$ws.Range("C9").Value2 = 1.406827509327035
$ws.Range("D9").Value2 = 2012
$ws.Range("E9").Value2 = 2.015050062499957

# row 10
$ws.Range("A10").Value2 = 40862
$ws.Range("B10").Value2 = 2011
$ws.Range("C10").Value2 = 0.899360810820804
$ws.Range("D10").Value2 = 2012
$ws.Range("E10").Value2 = 0.475544341751033

# row 11
$ws.Range("A11").Value2 = 41044
$ws.Range("B11").Value2 = 2012
$ws.Range("C11").Value2 = 1.153683074671208
$ws.Range("D11").Value2 = 2013
$ws.Range("E11").Value2 = 3.648892256099945

# row 12
$ws.Range("A12").Value2 = 41228
$ws.Range("B12").Value2 = 2012
$ws.Range("C12").Value2 = 0.9010266119894084
$ws.Range("D12").Value2 = 2013
$ws.Range("E12").Value2 = 1.506358095275817

# row 13
$ws.Range("A13").Value2 = 41409
$ws.Range("B13").Value2 = 2013
$ws.Range("C13").Value2 = 0.2186142574756467
$ws.Range("D13").Value2 = 2014
$ws.Range("E13").Value2 = 0.4006004000999708

# row 14
$ws.Range("A14").Value2 = 41592
$ws.Range("B14").Value2 = 2013
$ws.Range("C14").Value2 = 0.02019328874804938
$ws.Range("D14").Value2 = 2014
$ws.Range("E14").Value2 = -1.194807813319188

# row 15
$ws.Range("A15").Value2 = 41774
$ws.Range("B15").Value2 = 2014
$ws.Range("C15").Value2 = -0.8522658067264599
$ws.Range("D15").Value2 = 2015
$ws.Range("E15").Value2 = -3.551690943899999

# row 16
$ws.Range("A16").Value2 = 41957
$ws.Range("B16").Value2 = 2014
$ws.Range("C16").Value2 = 0.1729981757035093
$ws.Range("D16").Value2 = 2015
$ws.Range("E16").Value2 = 0.6265079396372775

# row 17
$ws.Range("A17").Value2 = 42137
$ws.Range("B17").Value2 = 2015
$ws.Range("C17").Value2 = -0.2262139320475365
$ws.Range("D17").Value2 = 2016
$ws.Range("E17").Value2 = -0.7976031983999876

# row 18
$ws.Range("A18").Value2 = 42321
$ws.Range("B18").Value2 = 2015
$ws.Range("C18").Value2 = 0.09752710595589686
$ws.Range("D18").Value2 = 2016
$ws.Range("E18").Value2 = -1.022506370243093

# row 19
$ws.Range("A19").Value2 = 42503
$ws.Range("B19").Value2 = 2016
$ws.Range("C19").Value2 = -0.6258176826215101
$ws.Range("D19").Value2 = 2017
$ws.Range("E19").Value2 = -0.3994003999000073

# row 20
$ws.Range("A20").Value2 = 42689
$ws.Range("B20").Value2 = 2016
$ws.Range("C20").Value2 = -0.5280591151586633
$ws.Range("D20").Value2 = 2017
$ws.Range("E20").Value2 = -0.7240982069264934

# row 21
$ws.Range("A21").Value2 = 42867
$ws.Range("B21").Value2 = 2017
$ws.Range("C21").Value2 = 0.3239252862367037
$ws.Range("D21").Value2 = 2018
$ws.Range("E21").Value2 = 1.609625625600009

# row 22
$ws.Range("A22").Value2 = 43053
$ws.Range("B22").Value2 = 2017
$ws.Range("C22").Value2 = 0.07201851318385799
$ws.Range("D22").Value2 = 2018
$ws.Range("E22").Value2 = 1.255028673974046

# row 23
$ws.Range("A23").Value2 = 43145
$ws.Range("B23").Value2 = 2018
$ws.Range("C23").Value2 = 1.456954732048321
$ws.Range("D23").Value2 = 2019
$ws.Range("E23").Value2 = 2.015050062499957

# row 24
$ws.Range("A24").Value2 = 43235
$ws.Range("B24").Value2 = 2018
$ws.Range("C24").Value2 = 0.5738128002843901
$ws.Range("D24").Value2 = 2019
$ws.Range("E24").Value2 = -0.3994003999000184

# row 25
$ws.Range("A25").Value2 = 43326
$ws.Range("B25").Value2 = 2018
$ws.Range("C25").Value2 = 0.3477859729380528
$ws.Range("D25").Value2 = 2019
$ws.Range("E25").Value2 = -1.516043567048941

# row 26
$ws.Range("A26").Value2 = 43418
$ws.Range("B26").Value2 = 2018
$ws.Range("C26").Value2 = 0.3727661260635617
$ws.Range("D26").Value2 = 2019
$ws.Range("E26").Value2 = -3.305525567352929

# row 27
$ws.Range("A27").Value2 = 43510
$ws.Range("B27").Value2 = 2019
$ws.Range("C27").Value2 = -0.7761690566734369
$ws.Range("D27").Value2 = 2020
$ws.Range("E27").Value2 = 0

# row 28
$ws.Range("A28").Value2 = 43600
$ws.Range("B28").Value2 = 2019
$ws.Range("C28").Value2 = -0.4781004700720293
$ws.Range("D28").Value2 = 2020
$ws.Range("E28").Value2 = 0.8024032015999882

# row 29
$ws.Range("A29").Value2 = 43691
$ws.Range("B29").Value2 = 2019
$ws.Range("C29").Value2 = -0.9254001004749823
$ws.Range("D29").Value2 = 2020
$ws.Range("E29").Value2 = -1.738778148048659

# row 30
$ws.Range("A30").Value2 = 43783
$ws.Range("B30").Value2 = 2019
$ws.Range("C30").Value2 = -0.801759526476209
$ws.Range("D30").Value2 = 2020
$ws.Range("E30").Value2 = 1.431264289671219

# row 31
$ws.Range("A31").Value2 = 43875
$ws.Range("B31").Value2 = 2020
$ws.Range("C31").Value2 = -1.172985875230903
$ws.Range("D31").Value2 = 2021
$ws.Range("E31").Value2 = -2.378486270400004

# row 32
$ws.Range("A32").Value2 = 43966
$ws.Range("B32").Value2 = 2020
$ws.Range("C32").Value2 = -1.197849743493773
$ws.Range("D32").Value2 = 2021
$ws.Range("E32").Value2 = -3.161804390399992

# row 33
$ws.Range("A33").Value2 = 44068
$ws.Range("B33").Value2 = 2020
$ws.Range("C33").Value2 = -1.503583188367719
$ws.Range("D33").Value2 = 2021
$ws.Range("E33").Value2 = 5.0514716327553

# row 34
$ws.Range("A34").Value2 = 44159
$ws.Range("B34").Value2 = 2020
$ws.Range("C34").Value2 = -1.103489789942047
$ws.Range("D34").Value2 = 2021
$ws.Range("E34").Value2 = 2.294626310579817

# row 35
$ws.Range("A35").Value2 = 44251
$ws.Range("B35").Value2 = 2021
$ws.Range("C35").Value2 = 3.668278063260222
$ws.Range("D35").Value2 = 2022
$ws.Range("E35").Value2 = 4.875032525328971

# row 36
$ws.Range("A36").Value2 = 44341
$ws.Range("B36").Value2 = 2021
$ws.Range("C36").Value2 = 1.064698711638945
$ws.Range("D36").Value2 = 2022
$ws.Range("E36").Value2 = -2.540956581357878

# row 37
$ws.Range("A37").Value2 = 44432
$ws.Range("B37").Value2 = 2021
$ws.Range("C37").Value2 = 1.067534122491809
$ws.Range("D37").Value2 = 2022
$ws.Range("E37").Value2 = 1.042084871410087

# row 38
$ws.Range("A38").Value2 = 44525
$ws.Range("B38").Value2 = 2021
$ws.Range("C38").Value2 = 0.9704846793491928
$ws.Range("D38").Value2 = 2022
$ws.Range("E38").Value2 = -0.902682013141165

# row 39
$ws.Range("A39").Value2 = 44617
$ws.Range("B39").Value2 = 2022
$ws.Range("C39").Value2 = -0.7181024432008964
$ws.Range("D39").Value2 = 2023
$ws.Range("E39").Value2 = -1.571815848026048

# row 40
$ws.Range("A40").Value2 = 44706
$ws.Range("B40").Value2 = 2022
$ws.Range("C40").Value2 = -1.524103236349472
$ws.Range("D40").Value2 = 2023
$ws.Range("E40").Value2 = -1.240907591477092

# row 41
$ws.Range("A41").Value2 = 44798
$ws.Range("B41").Value2 = 2022
$ws.Range("C41").Value2 = -0.9795431199870586
$ws.Range("D41").Value2 = 2023
$ws.Range("E41").Value2 = -0.5076503601560978

# row 42
$ws.Range("A42").Value2 = 44890
$ws.Range("B42").Value2 = 2022
$ws.Range("C42").Value2 = -0.7009264669202708
$ws.Range("D42").Value2 = 2023
$ws.Range("E42").Value2 = 1.220523709718857

# row 43
$ws.Range("A43").Value2 = 44981
$ws.Range("B43").Value2 = 2023
$ws.Range("C43").Value2 = 0.331635972555544
$ws.Range("D43").Value2 = 2024
$ws.Range("E43").Value2 = 1.49657333418427

# row 44
$ws.Range("A44").Value2 = 45071
$ws.Range("B44").Value2 = 2023
$ws.Range("C44").Value2 = 0.2854413827033664
$ws.Range("D44").Value2 = 2024
$ws.Range("E44").Value2 = -0.2470349027347551

# row 45
$ws.Range("A45").Value2 = 45163
$ws.Range("B45").Value2 = 2023
$ws.Range("C45").Value2 = 0.1829021030556488
$ws.Range("D45").Value2 = 2024
$ws.Range("E45").Value2 = -0.5273501419610804

# row 46
$ws.Range("A46").Value2 = 45254
$ws.Range("B46").Value2 = 2023
$ws.Range("C46").Value2 = 0.3928252664241905
$ws.Range("D46").Value2 = 2024
$ws.Range("E46").Value2 = 0.4517021897791018

# row 47
$ws.Range("A47").Value2 = 45345
$ws.Range("B47").Value2 = 2024
$ws.Range("C47").Value2 = 0.1776394553850924
$ws.Range("D47").Value2 = 2025
$ws.Range("E47").Value2 = -0.3176932480832284

# row 48
$ws.Range("A48").Value2 = 45436
$ws.Range("B48").Value2 = 2024
$ws.Range("C48").Value2 = 1.118108578853261
$ws.Range("D48").Value2 = 2025
$ws.Range("E48").Value2 = 1.532721825047534

# row 49
$ws.Range("A49").Value2 = 45534
$ws.Range("B49").Value2 = 2024
$ws.Range("C49").Value2 = 0.6979546684258597
$ws.Range("D49").Value2 = 2025
$ws.Range("E49").Value2 = -0.01999876157223746

# row 50
$ws.Range("A50").Value2 = 45618
$ws.Range("B50").Value2 = 2024
$ws.Range("C50").Value2 = 0.3224026462283813
$ws.Range("D50").Value2 = 2025
$ws.Range("E50").Value2 = -0.7185940249203049

# row 51
$ws.Range("A51").Value2 = 45713
$ws.Range("B51").Value2 = 2025
$ws.Range("C51").Value2 = -2.566037671324872
$ws.Range("D51").Value2 = 2026
$ws.Range("E51").Value2 = -1.090083898854388

# row 52
$ws.Range("A52").Value2 = 45800
$ws.Range("B52").Value2 = 2025
$ws.Range("C52").Value2 = -1.490505436658163
$ws.Range("D52").Value2 = 2026
$ws.Range("E52").Value2 = -0.3349088112516219

# row 53
$ws.Range("A53").Value2 = 45891
$ws.Range("B53").Value2 = 2025
$ws.Range("C53").Value2 = -2.11737366557071
$ws.Range("D53").Value2 = 2026
$ws.Range("E53").Value2 = -0.5919451648311758

